$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# The shared-string table dedupes identical text, so every cell that currently
# shares a given string must be updated together to keep that sharing (and end
# up producing the same single updated <si> entry in the saved workbook).

# "Latest HO Xliff Generate Date" (Overview col G) / "Correspond Handoff Datetime" (de-de col H)
# was "2016-08-17 00:14:22" -> "2016-08-17 00:15:12"
$wsOverview.Range("G3").Value = "2016-08-17 00:15:12"
$wsOverview.Range("G4").Value = "2016-08-17 00:15:12"
$wsDeDe.Range("H3").Value = "2016-08-17 00:15:12"
$wsDeDe.Range("H4").Value = "2016-08-17 00:15:12"

# "Status" (col E) was "ht" -> "mt" for both zh-cn and de-de
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn "Correspond Handoff Datetime" (col H) was "2016-08-17 00:14:18" -> "2016-08-17 00:15:02"
$wsZhCn.Range("H3").Value = "2016-08-17 00:15:02"
$wsZhCn.Range("H4").Value = "2016-08-17 00:15:02"

# zh-cn "Correspond Handback DateTime" (col K) was "2016-08-17 00:14:35" -> "2016-08-17 00:15:29"
$wsZhCn.Range("K3").Value = "2016-08-17 00:15:29"
$wsZhCn.Range("K4").Value = "2016-08-17 00:15:29"

# de-de "Correspond Handback DateTime" (col K) was "2016-08-17 00:14:42" -> "2016-08-17 00:15:36"
$wsDeDe.Range("K3").Value = "2016-08-17 00:15:36"
$wsDeDe.Range("K4").Value = "2016-08-17 00:15:36"
